$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new columns: I = banEmailUpdate, H = banThirdParty
# (shared-string table order follows the order cells are first written,
# so write column I's strings before column H's to match the target order)
$ws.Range("I1").Value = "banEmailUpdate"
$ws.Range("I2").Value = "选填，0或1，1表示禁止更改邮箱"

$ws.Range("H1").Value = "banThirdParty"
$ws.Range("H2").Value = "选填，0或1，1表示禁止使用第三方登录"

# Move the active selection like the author left it
$ws.Range("K12").Select()
